# edit.ps1
# Applies the update described by the commit diff:
#  1. Column C ("Förändrad") for all data rows (2-82) changes from 46070 to 46072.
#  2. Rows 31-61 are re-ordered (the underlying list was re-sorted/updated), which
#     changes the Beteckning (A), Datum (B), Markägare (F) and Area (G) values that
#     live on each row number. We reproduce the resulting per-row values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Column C: 46070 -> 46072 for every data row (2 through 82) ---
$ws.Range("C2:C82").Value = 46072

# --- 2. Rows 31-61: new contents after the reorder/update ---

# Row 31
$ws.Range("A31").Value = "A 21112-2025"
$ws.Range("B31").Value = 45777
$ws.Range("G31").Value = 7.2

# Row 32
$ws.Range("A32").Value = "A 18698-2025"
$ws.Range("B32").Value = 45763
$ws.Range("G32").Value = 29.7

# Row 33
$ws.Range("A33").Value = "A 61066-2024"
$ws.Range("B33").Value = 45645.49744212963
$ws.Range("G33").Value = 4.5

# Row 34 (unchanged besides column C, already applied above)

# Row 35
$ws.Range("A35").Value = "A 50590-2023"
$ws.Range("B35").Value = 45217.40167824074
$ws.Range("G35").Value = 1.3

# Row 36
$ws.Range("A36").Value = "A 62460-2022"
$ws.Range("B36").Value = 44916
$ws.Range("G36").Value = 13.4

# Row 37 (unchanged besides column C, already applied above)

# Row 38
$ws.Range("A38").Value = "A 33183-2025"
$ws.Range("B38").Value = 45840.59206018518
$ws.Range("G38").Value = 2.8

# Row 39
$ws.Range("A39").Value = "A 34663-2025"
$ws.Range("B39").Value = 45848.37454861111
$ws.Range("G39").Value = 5.1

# Row 40 (gains a Markägare value in column F)
$ws.Range("A40").Value = "A 50338-2025"
$ws.Range("B40").Value = 45944.52140046296
$ws.Range("F40").Value = "SCA"
$ws.Range("G40").Value = 22.6

# Row 41
$ws.Range("A41").Value = "A 34664-2025"
$ws.Range("B41").Value = 45848.37857638889
$ws.Range("G41").Value = 8.4

# Row 42
$ws.Range("A42").Value = "A 5533-2024"
$ws.Range("B42").Value = 45334
$ws.Range("G42").Value = 2.4

# Row 43 (loses its Markägare value in column F)
$ws.Range("A43").Value = "A 9605-2024"
$ws.Range("B43").Value = 45361.44230324074
$ws.Range("F43").ClearContents()
$ws.Range("G43").Value = 8.4

# Row 44
$ws.Range("A44").Value = "A 5411-2024"
$ws.Range("B44").Value = 45331
$ws.Range("G44").Value = 3.9

# Row 45
$ws.Range("A45").Value = "A 28707-2022"
$ws.Range("B45").Value = 44748
$ws.Range("G45").Value = 15.2

# Row 46
$ws.Range("A46").Value = "A 20893-2025"
$ws.Range("B46").Value = 45776
$ws.Range("G46").Value = 27.3

# Row 47
$ws.Range("A47").Value = "A 21184-2025"
$ws.Range("B47").Value = 45777
$ws.Range("G47").Value = 20.1

# Row 48 (unchanged besides column C, already applied above; keeps "SCA" in F)

# Row 49 (unchanged besides column C, already applied above)

# Row 50
$ws.Range("A50").Value = "A 2664-2026"
$ws.Range("B50").Value = 46037.64078703704
$ws.Range("G50").Value = 31.3

# Row 51
$ws.Range("A51").Value = "A 39271-2023"
$ws.Range("B51").Value = 45162
$ws.Range("G51").Value = 51.9

# Row 52
$ws.Range("A52").Value = "A 2658-2026"
$ws.Range("B52").Value = 46037.63512731482
$ws.Range("G52").Value = 22.5

# Row 53
$ws.Range("A53").Value = "A 18712-2025"
$ws.Range("B53").Value = 45763
$ws.Range("G53").Value = 66.5

# Row 54 (gains a Markägare value in column F)
$ws.Range("A54").Value = "A 56168-2023"
$ws.Range("B54").Value = 45236
$ws.Range("F54").Value = "Sveaskog"
$ws.Range("G54").Value = 6.6

# Row 55
$ws.Range("A55").Value = "A 22494-2024"
$ws.Range("B55").Value = 45447
$ws.Range("G55").Value = 6

# Row 56 (loses its Markägare value in column F)
$ws.Range("A56").Value = "A 49704-2024"
$ws.Range("B56").Value = 45596
$ws.Range("F56").ClearContents()
$ws.Range("G56").Value = 17.4

# Row 57
$ws.Range("A57").Value = "A 22491-2024"
$ws.Range("B57").Value = 45447
$ws.Range("G57").Value = 9

# Row 58
$ws.Range("A58").Value = "A 16024-2024"
$ws.Range("B58").Value = 45405
$ws.Range("G58").Value = 2.3

# Row 59
$ws.Range("A59").Value = "A 2660-2026"
$ws.Range("B59").Value = 46037
$ws.Range("G59").Value = 1.1

# Row 60 (gains a Markägare value in column F)
$ws.Range("A60").Value = "A 47730-2024"
$ws.Range("B60").Value = 45588
$ws.Range("F60").Value = "Allmännings- och besparingsskogar"
$ws.Range("G60").Value = 3.7

# Row 61 (loses its Markägare value in column F)
$ws.Range("A61").Value = "A 148-2026"
$ws.Range("B61").Value = 46024.6446875
$ws.Range("F61").ClearContents()
$ws.Range("G61").Value = 10.2
